$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($row1, $row2)

    $range1 = $ws.Range("B$row1`:AB$row1")
    $range2 = $ws.Range("B$row2`:AB$row2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}

Swap-Rows 38 39
Swap-Rows 84 85
Swap-Rows 246 247
